$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1, copying the existing header style (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New "Save" data column values (H2:H6)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
